$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LIST")

$ws.Range("A3").Value = "RT.MAT.001.REC"
$ws.Range("A4").Value = "RT.MAT.001.SUP"

$ws.Range("C11").Select()
